# Finished analysis for preliminary doubling time dilutions:
# swap the GFP 1x / MC 1x labels between the alternating dilution rows,
# clear out the now-unused row labels (E,F,G,H) and their row contents,
# and leave the selection on the last edited range.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 & 4 (dilution rows "A" and "C") were labeled GFP 1x -> now MC 1x
$ws.Range("B2:J2").Value = "MC 1x"
$ws.Range("B4:J4").Value = "MC 1x"

# Rows 3 & 5 (dilution rows "B" and "D") were labeled MC 1x -> now GFP 1x
$ws.Range("B3:J3").Value = "GFP 1x"
$ws.Range("B5:J5").Value = "GFP 1x"

# Remove the row labels (E, F, G, H) for the unused rows 6-9
$ws.Range("A6:A9").ClearContents()

# Fully clear out the contents/formatting of the unused row data in rows 7 & 8
$ws.Range("B7:J7").Clear()
$ws.Range("B8:J8").Clear()

# Leave the final selection on the last-touched range
$ws.Range("B5:J5").Select()
